$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row = 12; I = 'sd'; J = 'Statement-non-opinion' },
    @{ Row = 19; I = 'sd'; J = 'Statement-non-opinion' },
    @{ Row = 20; I = 'b'; J = 'Acknowledge (Backchannel)' },
    @{ Row = 22; I = 'ba'; J = 'Appreciation' },
    @{ Row = 34; I = 'sd'; J = 'Statement-non-opinion' },
    @{ Row = 37; I = 'sd'; J = 'Statement-non-opinion' },
    @{ Row = 38; I = 'sd'; J = 'Statement-non-opinion' },
    @{ Row = 46; I = 'sd'; J = 'Statement-non-opinion' },
    @{ Row = 51; I = 'sd'; J = 'Statement-non-opinion' },
    @{ Row = 53; I = 'ba'; J = 'Appreciation' },
    @{ Row = 54; I = 'sd'; J = 'Statement-non-opinion' },
    @{ Row = 57; I = 'sd'; J = 'Statement-non-opinion' },
    @{ Row = 59; I = 'aa'; J = 'Agree/Accept' },
    @{ Row = 77; I = 'sv'; J = 'Statement-opinion' },
    @{ Row = 78; I = 'sd'; J = 'Statement-non-opinion' },
    @{ Row = 84; I = 'qy'; J = 'Yes-No-Question' },
    @{ Row = 109; I = 'ba'; J = 'Appreciation' },
    @{ Row = 112; I = 'sv'; J = 'Statement-opinion' },
    @{ Row = 125; I = '%'; J = 'Uninterpretable' },
    @{ Row = 129; I = 'sd'; J = 'Statement-non-opinion' },
    @{ Row = 133; I = 'sd'; J = 'Statement-non-opinion' },
    @{ Row = 137; I = 'sd'; J = 'Statement-non-opinion' },
    @{ Row = 140; I = 'aa'; J = 'Agree/Accept' },
    @{ Row = 141; I = 'sv'; J = 'Statement-opinion' },
    @{ Row = 156; I = '%'; J = 'Uninterpretable' },
    @{ Row = 163; I = 'sd'; J = 'Statement-non-opinion' },
    @{ Row = 165; I = 'sd'; J = 'Statement-non-opinion' },
    @{ Row = 170; I = 'sv'; J = 'Statement-opinion' },
    @{ Row = 174; I = 'sd'; J = 'Statement-non-opinion' },
    @{ Row = 179; I = 'sv'; J = 'Statement-opinion' },
    @{ Row = 189; I = 'sd'; J = 'Statement-non-opinion' },
    @{ Row = 196; I = 'ba'; J = 'Appreciation' },
    @{ Row = 202; I = 'sv'; J = 'Statement-opinion' },
    @{ Row = 204; I = 'sd'; J = 'Statement-non-opinion' },
    @{ Row = 213; I = 'sd'; J = 'Statement-non-opinion' },
    @{ Row = 216; I = 'sd'; J = 'Statement-non-opinion' },
    @{ Row = 221; I = 'sd'; J = 'Statement-non-opinion' },
    @{ Row = 233; I = 'sd'; J = 'Statement-non-opinion' },
    @{ Row = 272; I = 'aa'; J = 'Agree/Accept' },
    @{ Row = 275; I = 'aa'; J = 'Agree/Accept' },
    @{ Row = 284; I = 'sd'; J = 'Statement-non-opinion' },
    @{ Row = 295; I = 'sv'; J = 'Statement-opinion' },
    @{ Row = 296; I = 'sd'; J = 'Statement-non-opinion' },
    @{ Row = 303; I = 'sd'; J = 'Statement-non-opinion' },
    @{ Row = 306; I = 'sv'; J = 'Statement-opinion' },
    @{ Row = 311; I = 'sv'; J = 'Statement-opinion' },
    @{ Row = 328; I = 'sd'; J = 'Statement-non-opinion' },
    @{ Row = 336; I = 'sd'; J = 'Statement-non-opinion' },
    @{ Row = 341; I = 'b'; J = 'Acknowledge (Backchannel)' },
    @{ Row = 343; I = 'sd'; J = 'Statement-non-opinion' },
    @{ Row = 345; I = 'sd'; J = 'Statement-non-opinion' },
    @{ Row = 370; I = 'aa'; J = 'Agree/Accept' },
    @{ Row = 372; I = 'sv'; J = 'Statement-opinion' },
    @{ Row = 374; I = 'aa'; J = 'Agree/Accept' },
    @{ Row = 385; I = 'sd'; J = 'Statement-non-opinion' },
    @{ Row = 388; I = 'sd'; J = 'Statement-non-opinion' },
    @{ Row = 392; I = 'sd'; J = 'Statement-non-opinion' },
    @{ Row = 395; I = 'sd'; J = 'Statement-non-opinion' },
    @{ Row = 429; I = 'sd'; J = 'Statement-non-opinion' },
    @{ Row = 439; I = '%'; J = 'Uninterpretable' },
    @{ Row = 442; I = '%'; J = 'Uninterpretable' },
    @{ Row = 446; I = 'aa'; J = 'Agree/Accept' },
    @{ Row = 448; I = 'aa'; J = 'Agree/Accept' },
    @{ Row = 452; I = 'sd'; J = 'Statement-non-opinion' },
    @{ Row = 460; I = 'sd'; J = 'Statement-non-opinion' },
    @{ Row = 472; I = 'sd'; J = 'Statement-non-opinion' },
    @{ Row = 480; I = 'aa'; J = 'Agree/Accept' },
    @{ Row = 488; I = 'b'; J = 'Acknowledge (Backchannel)' },
    @{ Row = 490; I = 'b'; J = 'Acknowledge (Backchannel)' },
    @{ Row = 499; I = 'aa'; J = 'Agree/Accept' },
    @{ Row = 522; I = 'sd'; J = 'Statement-non-opinion' }
)

foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, 9).Value = $u.I
    $ws.Cells.Item($u.Row, 10).Value = $u.J
}
